$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the conditional formatting that used to live on column G (Status)
# before that column goes away.
$ws.Range("G1:G1047927").FormatConditions.Delete()

# Reorder the sheet so the column we are keeping ("Name", originally D)
# sits right after "Site", then drop the now-trailing "Comment" / "Title" /
# "Template" / "Status" columns in a single delete.
$ws.Columns("D:D").Cut() | Out-Null
$ws.Columns("C:C").Insert() | Out-Null
$ws.Columns("D:G").Delete()

# Shrink the table ("Tabelle4") down to the three remaining columns and
# make sure its header picks up the "Name" caption that now lives in C1.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C1047927"))
$ws.Range("C1").Value = "Name"

# Move the active selection like the saved workbook shows.
$ws.Range("I20").Select() | Out-Null
